$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1956521739130435
$ws.Range("C2").Value = 0.5543478260869565
$ws.Range("J2").Value = 0.007246376811594203
$ws.Range("P2").Value = 0.1268115942028986
$ws.Range("S2").Value = 0.1159420289855072
$ws.Range("B3").Value = 0.02531645569620253
$ws.Range("C3").Value = 0.0189873417721519
$ws.Range("J3").Value = 0.01265822784810127
$ws.Range("P3").Value = 0.759493670886076
$ws.Range("S3").Value = 0.1835443037974684
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2619047619047619
$ws.Range("B6").Value = 0.1064814814814815
$ws.Range("D6").Value = 0.01388888888888889
$ws.Range("F6").Value = 0.06018518518518518
$ws.Range("J6").Value = 0.1990740740740741
$ws.Range("O6").Value = 0.01851851851851852
$ws.Range("Q6").Value = 0.125
$ws.Range("R6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.3935185185185185
$ws.Range("B7").Value = 0.1057692307692308
$ws.Range("D7").Value = 0.01923076923076923
$ws.Range("F7").Value = 0.05288461538461538
$ws.Range("J7").Value = 0.1346153846153846
$ws.Range("O7").Value = 0.01442307692307692
$ws.Range("Q7").Value = 0.2163461538461539
$ws.Range("R7").Value = 0.07211538461538461
$ws.Range("S7").Value = 0.3846153846153846
$ws.Range("B8").Value = 0.06762295081967214
$ws.Range("D8").Value = 0.02459016393442623
$ws.Range("F8").Value = 0.07172131147540983
$ws.Range("J8").Value = 0.1331967213114754
$ws.Range("O8").Value = 0.01639344262295082
$ws.Range("Q8").Value = 0.1782786885245902
$ws.Range("R8").Value = 0.07172131147540983
$ws.Range("S8").Value = 0.4364754098360656
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.01578947368421053
$ws.Range("F9").Value = 0.06315789473684211
$ws.Range("J9").Value = 0.1578947368421053
$ws.Range("O9").Value = 0.01578947368421053
$ws.Range("Q9").Value = 0.1736842105263158
$ws.Range("R9").Value = 0.05263157894736842
$ws.Range("S9").Value = 0.4210526315789473
$ws.Range("B10").Value = 0.1002465078060805
$ws.Range("D10").Value = 0.0180772391125719
$ws.Range("E10").Value = 0.0008216926869350862
$ws.Range("F10").Value = 0.0714872637633525
$ws.Range("J10").Value = 0.1273623664749384
$ws.Range("O10").Value = 0.01725554642563681
$ws.Range("Q10").Value = 0.2292522596548891
$ws.Range("R10").Value = 0.06820049301561217
$ws.Range("S10").Value = 0.3672966310599836
$ws.Range("F11").Value = 0.003257328990228013
$ws.Range("G11").Value = 0.1693811074918567
$ws.Range("J11").Value = 0.07491856677524431
$ws.Range("K11").Value = 0.2280130293159609
$ws.Range("L11").Value = 0.5146579804560261
$ws.Range("S11").Value = 0.009771986970684038
$ws.Range("G12").Value = 0.7559523809523809
$ws.Range("J12").Value = 0.1845238095238095
$ws.Range("L12").Value = 0.04761904761904762
$ws.Range("S12").Value = 0.0119047619047619
$ws.Range("G13").Value = 0.8409090909090909
$ws.Range("J13").Value = 0.09090909090909091
$ws.Range("S13").Value = 0.06818181818181818
$ws.Range("F15").Value = 0.01762114537444934
$ws.Range("H15").Value = 0.1497797356828194
$ws.Range("I15").Value = 0.04845814977973568
$ws.Range("J15").Value = 0.3348017621145374
$ws.Range("K15").Value = 0.05726872246696035
$ws.Range("M15").Value = 0.00881057268722467
$ws.Range("O15").Value = 0.08370044052863436
$ws.Range("S15").Value = 0.2995594713656388
$ws.Range("F16").Value = 0.01685393258426966
$ws.Range("H16").Value = 0.1853932584269663
$ws.Range("I16").Value = 0.07865168539325842
$ws.Range("J16").Value = 0.398876404494382
$ws.Range("K16").Value = 0.0898876404494382
$ws.Range("M16").Value = 0.01685393258426966
$ws.Range("O16").Value = 0.05617977528089887
$ws.Range("S16").Value = 0.1573033707865168
$ws.Range("F17").Value = 0.01720430107526882
$ws.Range("H17").Value = 0.1763440860215054
$ws.Range("I17").Value = 0.08387096774193549
$ws.Range("J17").Value = 0.3978494623655914
$ws.Range("K17").Value = 0.09677419354838709
$ws.Range("M17").Value = 0.01935483870967742
$ws.Range("O17").Value = 0.05591397849462366
$ws.Range("S17").Value = 0.1526881720430107
$ws.Range("F18").Value = 0.01226993865030675
$ws.Range("H18").Value = 0.2269938650306748
$ws.Range("I18").Value = 0.09202453987730061
$ws.Range("J18").Value = 0.3803680981595092
$ws.Range("K18").Value = 0.09815950920245399
$ws.Range("M18").Value = 0.01226993865030675
$ws.Range("O18").Value = 0.03067484662576687
$ws.Range("S18").Value = 0.147239263803681
$ws.Range("F19").Value = 0.008759124087591242
$ws.Range("H19").Value = 0.2211678832116788
$ws.Range("I19").Value = 0.08102189781021898
$ws.Range("J19").Value = 0.327007299270073
$ws.Range("K19").Value = 0.1072992700729927
$ws.Range("M19").Value = 0.02116788321167883
$ws.Range("N19").Value = 0.00072992700729927
$ws.Range("O19").Value = 0.07591240875912408
$ws.Range("S19").Value = 0.1569343065693431
